$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.575.61'
$ws.Cells.Item(2, 5).Value = '  +0.44%  '

$ws.Cells.Item(3, 4).Value = '3.444.26'
$ws.Cells.Item(3, 5).Value = '  +2.75%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '579.46'
$ws.Cells.Item(5, 5).Value = '  +1.41%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '147.11'
$ws.Cells.Item(6, 5).Value = '  +7.93%  '

$ws.Cells.Item(7, 4).Value = '3.444.89'
$ws.Cells.Item(7, 5).Value = '  +2.81%  '

$ws.Cells.Item(8, 5).Value = '  +0.02%  '

$ws.Cells.Item(9, 5).Value = '  +0.68%  '

$ws.Cells.Item(10, 5).Value = '  +2.50%  '

$ws.Cells.Item(11, 5).Value = '  +0.25%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.388'
$ws.Cells.Item(12, 5).Value = '  -0.14%  '

$ws.Cells.Item(13, 4).Value = '4.031.62'
$ws.Cells.Item(13, 5).Value = '  +2.62%  '

$ws.Cells.Item(14, 5).Value = '  -0.84%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '27.77'
$ws.Cells.Item(15, 5).Value = '  +7.56%  '

$ws.Cells.Item(16, 5).Value = '  +0.14%  '

$ws.Cells.Item(17, 4).Value = '3.444.82'

$ws.Cells.Item(18, 4).Value = '61.680.87'
$ws.Cells.Item(18, 5).Value = '  +0.37%  '

$ws.Cells.Item(19, 5).Value = '  +8.20%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.08'
$ws.Cells.Item(20, 5).Value = '  +1.19%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '9.48'
$ws.Cells.Item(21, 5).Value = '  +2.16%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '387.03'
$ws.Cells.Item(22, 5).Value = '  +3.25%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.564'
$ws.Cells.Item(23, 5).Value = '  +2.49%  '

$ws.Cells.Item(24, 4).Value = '3.588.93'

$ws.Cells.Item(25, 5).Value = '  +0.09%  '

$ws.Cells.Item(26, 2).Value = 'Litecoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '72.28'
$ws.Cells.Item(26, 5).Value = '  +2.00%  '

$ws.Cells.Item(27, 2).Value = 'LEO'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '5.74'
$ws.Cells.Item(27, 5).Value = '  +0.03%  '

$ws.Cells.Item(28, 5).Value = '  -0.44%  '

$ws.Cells.Item(29, 2).Value = 'RenderToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.85'
$ws.Cells.Item(29, 5).Value = '  +5.09%  '

$ws.Cells.Item(30, 2).Value = 'Kaspa'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.175'
$ws.Cells.Item(30, 5).Value = '  +7.14%  '

$ws.Cells.Item(31, 5).Value = '  -10.82%  '

$ws.Cells.Item(32, 5).Value = '  -0.10%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '8.20'
$ws.Cells.Item(33, 5).Value = '  +0.80%  '

$ws.Cells.Item(34, 5).Value = '  +1.26%  '

$ws.Cells.Item(35, 5).Value = '  -0.03%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '24.15'
$ws.Cells.Item(36, 5).Value = '  +2.81%  '

$ws.Cells.Item(37, 5).Value = '  +1.05%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '7.01'
$ws.Cells.Item(38, 5).Value = '  +3.33%  '

$ws.Cells.Item(39, 5).Value = '  +2.34%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '166.11'
$ws.Cells.Item(40, 5).Value = '  +0.83%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0789'
$ws.Cells.Item(41, 5).Value = '  +3.27%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '25.88'
$ws.Cells.Item(42, 5).Value = '  +7.60%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.790'
$ws.Cells.Item(43, 5).Value = '  +2.72%  '

$ws.Cells.Item(44, 5).Value = '  -0.15%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '4.49'
$ws.Cells.Item(45, 5).Value = '  +2.98%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.73'
$ws.Cells.Item(46, 5).Value = '  +0.66%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '42.01'
$ws.Cells.Item(47, 5).Value = '  +1.71%  '

$ws.Cells.Item(48, 4).Value = '2.629.24'
$ws.Cells.Item(48, 5).Value = '  +11.25%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.17'
$ws.Cells.Item(49, 5).Value = '  -2.33%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '23.79'
$ws.Cells.Item(50, 5).Value = '  +4.72%  '

$ws.Cells.Item(51, 5).Value = '  +0.39%  '
